$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 previously held the "-3" message id entry. init_msg_db no longer opens
# a duplicate copy of the database, so row 4 now reuses the id=1 message
# (same short description as row 2) while keeping its own long
# description/remedy text, and the numeric id becomes a plain 3.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "short db message id 1,"

# The now-unused "short db message id -3," shared string is dropped
# automatically when the workbook is saved since no cell references it
# anymore.

# Reflect the new selection left behind in the sheet (B2:B4, active cell B2).
$ws.Range("B2:B4").Select() | Out-Null
